{"js": "// Replace the date and each division problem's text with its new value.\n// All old values are unique within the document, so a direct search +\n// Replace for each pair is safe and preserves run formatting (rFonts, sz).\nconst replacements = [\n  [\"2025-06-01 Sunday\", \"2025-06-02 Monday\"],\n  [\"848\u00f76=\", \"825\u00f73=\"],\n  [\"681\u00f73=\", \"258\u00f78=\"],\n  [\"293\u00f77=\", \"205\u00f78=\"],\n  [\"187\u00f78=\", \"981\u00f76=\"],\n  [\"369\u00f74=\", \"318\u00f73=\"],\n  [\"861\u00f79=\", \"836\u00f74=\"],\n  [\"558\u00f78=\", \"432\u00f76=\"],\n  [\"561\u00f78=\", \"908\u00f73=\"],\n  [\"581\u00f78=\", \"625\u00f74=\"],\n  [\"846\u00f73=\", \"676\u00f73=\"],\n  [\"803\u00f78=\", \"190\u00f79=\"],\n  [\"235\u00f74=\", \"881\u00f77=\"],\n  [\"904\u00f72=\", \"937\u00f75=\"],\n  [\"343\u00f79=\", \"294\u00f72=\"],\n  [\"190\u00f75=\", \"143\u00f77=\"],\n  [\"852\u00f78=\", \"842\u00f79=\"],\n  [\"289\u00f77=\", \"741\u00f74=\"],\n  [\"194\u00f73=\", \"590\u00f74=\"],\n  [\"129\u00f73=\", \"106\u00f74=\"],\n  [\"488\u00f78=\", \"986\u00f79=\"],\n  [\"903\u00f74=\", \"757\u00f76=\"],\n  [\"123\u00f79=\", \"522\u00f79=\"],\n  [\"534\u00f78=\", \"137\u00f77=\"],\n  [\"468\u00f74=\", \"518\u00f78=\"],\n  [\"202\u00f76=\", \"101\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each division problem's text with its new value\n# using Find/Replace over the whole document. Every old value is unique\n# in the document, so a single wdReplaceAll pass per pair is safe and\n# keeps each run's original formatting (rFonts, sz) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-01 Sunday\", \"2025-06-02 Monday\"),\n    @(\"848\u00f76=\", \"825\u00f73=\"),\n    @(\"681\u00f73=\", \"258\u00f78=\"),\n    @(\"293\u00f77=\", \"205\u00f78=\"),\n    @(\"187\u00f78=\", \"981\u00f76=\"),\n    @(\"369\u00f74=\", \"318\u00f73=\"),\n    @(\"861\u00f79=\", \"836\u00f74=\"),\n    @(\"558\u00f78=\", \"432\u00f76=\"),\n    @(\"561\u00f78=\", \"908\u00f73=\"),\n    @(\"581\u00f78=\", \"625\u00f74=\"),\n    @(\"846\u00f73=\", \"676\u00f73=\"),\n    @(\"803\u00f78=\", \"190\u00f79=\"),\n    @(\"235\u00f74=\", \"881\u00f77=\"),\n    @(\"904\u00f72=\", \"937\u00f75=\"),\n    @(\"343\u00f79=\", \"294\u00f72=\"),\n    @(\"190\u00f75=\", \"143\u00f77=\"),\n    @(\"852\u00f78=\", \"842\u00f79=\"),\n    @(\"289\u00f77=\", \"741\u00f74=\"),\n    @(\"194\u00f73=\", \"590\u00f74=\"),\n    @(\"129\u00f73=\", \"106\u00f74=\"),\n    @(\"488\u00f78=\", \"986\u00f79=\"),\n    @(\"903\u00f74=\", \"757\u00f76=\"),\n    @(\"123\u00f79=\", \"522\u00f79=\"),\n    @(\"534\u00f78=\", \"137\u00f77=\"),\n    @(\"468\u00f74=\", \"518\u00f78=\"),\n    @(\"202\u00f76=\", \"101\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
